$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row: "4.simple_e_p_mix_random_test" with its (wrapped) description.
$ws.Range("A5").Value = "4.simple_e_p_mix_random_test"
$ws.Range("B5").Value = ">express packet and preemptable packet random `n>packet length of express  and  preemptable packet are both random`n"

# Widen the description column and wrap the new cell's text, growing the row to fit.
# (47.0 is what lands the saved OOXML <col> width nearest the authored 47.75
# once the host's character/pixel grid rounding is applied on save.)
$ws.Columns.Item(2).ColumnWidth = 47.0
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 57

# Leave the selection where the author left it when finishing the edit.
[void]$ws.Range("B8").Select()
